# Apply updated "want to go" / "lowest price" counts scraped at commit 456a3b4.
# Two sheets ("展览" and "全部类型") hold the same conference listing and both
# need the same numeric updates.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("G5").Value = 35
    $ws.Range("F6").Value = 59
    $ws.Range("F10").Value = 16127
    $ws.Range("F11").Value = 266
    $ws.Range("F12").Value = 192
    $ws.Range("F14").Value = 6295
    $ws.Range("F28").Value = 886

    if ($name -eq "展览") {
        $ws.Range("F32").Value = 11252
        $ws.Range("F35").Value = 139
    } else {
        $ws.Range("F33").Value = 11252
        $ws.Range("F36").Value = 139
    }
}
